# The document has a paragraph (highlighted yellow) whose only content is an
# inline picture (Picture 19 / rId11). The edit removes that picture run,
# leaving the paragraph mark (and its formatting) behind -- exactly like a
# user selecting the picture in Word and pressing Delete.
$d = $word.ActiveDocument

# Walk the InlineShapes collection backwards and drop every inline picture;
# this mirrors "select picture -> Delete" without hard-coding a single index.
for ($i = $d.InlineShapes.Count; $i -ge 1; $i--) {
    $shape = $d.InlineShapes($i)
    $shape.Delete()
}
